# "Updated SFI GM,VR and UF"
#
# A2 held a hyperlink to the old RThree form URL (formid=220136). The new
# commit points it at a new form (formid=829192) and the cell is no longer a
# hyperlink (its "Hyperlink" style is dropped back to Normal). Three new
# lookup cells are also added to row 2: K2/AU2 ("formshow_06112023") and
# AR2 ("voice_record_06112023").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: drop the hyperlink, point the text at the new form URL, and reset the
# cell's style back to Normal (it no longer carries the Hyperlink format).
$ws.Range("A2").Hyperlinks.Delete()
$ws.Range("A2").Value = "https://rthree.live/showform?formid=829192&nurams=bot1"
$ws.Range("A2").Style = "Normal"

# New values appended to row 2 (voice_record_06112023 is the first of the
# two new strings to be introduced, so it keeps the earlier shared-string
# slot; formshow_06112023 -- used by both K2 and AU2 -- comes after it).
$ws.Range("AR2").Value = "voice_record_06112023"
$ws.Range("K2").Value = "formshow_06112023"
$ws.Range("AU2").Value = "formshow_06112023"

# The sheet's active selection moved from G1 to A2.
$ws.Range("A2").Select()
